# Open Source Task workbook: add a submission row (name / email / repo link)
# with the email and repo link rendered as hyperlinks, matching the
# author's "edit in xlsx file by ahmed qabeel" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 values + hyperlinks (order matches the shared-string table
#     produced by the original authoring session: repo link, name, email)
$ws.Range("C2").Value = "https://github.com/dohaqabeel/Security-Task.git"
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/dohaqabeel/Security-Task.git") | Out-Null

$ws.Range("A2").Value = "أحمد ممدوح أمين قابيل"

$ws.Range("B2").Value = "doo7a008@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:doo7a008@gmail.com") | Out-Null

# --- Selection, matching the saved view --------------------------------
$ws.Range("C5").Select() | Out-Null
